$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.2564814264887687
$ws.Range("A2").Value = -0.0059999999671198623
$ws.Range("A3").Value = -0.0039999999707198697
$ws.Range("A4").Value = -0.007999999946630254
$ws.Range("A5").Value = -0.0029999999695942137
$ws.Range("A6").Value = -0.0019999999672659641
$ws.Range("A7").Value = -0.0099999999257756045
$ws.Range("A8").Value = -0.0099999999261357608
$ws.Range("A9").Value = -0.0019999999694699788
$ws.Range("A10").Value = 0.03688940217286607
$ws.Range("A11").Value = -0.0029999999659882093
$ws.Range("A12").Value = -0.0034999999629192757
$ws.Range("A13").Value = -0.0034999999615248356
$ws.Range("A14").Value = -0.0079999999377626807
$ws.Range("A15").Value = -0.00099999997434174759
$ws.Range("A16").Value = -0.0019999999691386883
$ws.Range("A17").Value = -0.0019999999693691706
$ws.Range("A18").Value = 0.0025005020747723705
$ws.Range("A19").Value = -0.0039999999766973104
$ws.Range("A20").Value = -0.066578075142585291
$ws.Range("A21").Value = -0.0039999999666040509
$ws.Range("A22").Value = -0.0039999999664335206
$ws.Range("A23").Value = -0.0049999999638661308
$ws.Range("A24").Value = -0.019999999882597486
$ws.Range("A25").Value = -0.019999999881148867
$ws.Range("A26").Value = -0.002499999963253785
$ws.Range("A27").Value = -0.0024999999612012047
$ws.Range("A28").Value = -0.0019999999547604119
$ws.Range("A29").Value = -0.0069999999223613329
$ws.Range("A30").Value = -0.059999999641908275
$ws.Range("A31").Value = -0.006999999916965649
$ws.Range("A32").Value = -0.009999999900943024
$ws.Range("A33").Value = -0.0039999999322386515
